# Updates D (Price) and E (Volume(1h)) columns on the crypto tracker sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern for D-column (Price) cells: some new values are plain
# numeric-looking strings (e.g. "592.92"); Excel would auto-convert those
# to numbers on assignment. Briefly force the cell to Text ("@") so the
# literal string is preserved, then ClearFormats() so no extra number
# format / style survives on the cell (matching the original unstyled cells).
function Set-TextPrice($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextPrice "D2" "67.331.87"
$ws.Range("E2").Value = "  -3.08%  "

Set-TextPrice "D3" "3.766.23"
$ws.Range("E3").Value = "  -1.26%  "

$ws.Range("E4").Value = "  +0.00%  "

Set-TextPrice "D5" "592.92"
$ws.Range("E5").Value = "  -1.15%  "

Set-TextPrice "D6" "165.87"
$ws.Range("E6").Value = "  -3.46%  "

Set-TextPrice "D7" "3.764.87"
$ws.Range("E7").Value = "  -1.28%  "

$ws.Range("E9").Value = "  -1.85%  "

$ws.Range("E10").Value = "  -3.00%  "

$ws.Range("E11").Value = "  -2.63%  "

$ws.Range("E13").Value = "  -4.30%  "

Set-TextPrice "D14" "35.84"
$ws.Range("E14").Value = "  -3.05%  "

Set-TextPrice "D15" "4.398.71"
$ws.Range("E15").Value = "  -1.21%  "

Set-TextPrice "D16" "3.751.96"
$ws.Range("E16").Value = "  -1.55%  "

Set-TextPrice "D17" "67.326.06"
$ws.Range("E17").Value = "  -3.00%  "

Set-TextPrice "D18" "17.76"
$ws.Range("E18").Value = "  -3.03%  "

$ws.Range("E19").Value = "  -0.08%  "

$ws.Range("E20").Value = "  -2.48%  "

Set-TextPrice "D21" "10.21"
$ws.Range("E21").Value = "  -8.59%  "

Set-TextPrice "D22" "456.89"
$ws.Range("E22").Value = "  -3.71%  "

Set-TextPrice "D23" "0.697"
$ws.Range("E23").Value = "  -1.98%  "

$ws.Range("E24").Value = "  +0.96%  "

Set-TextPrice "D25" "83.13"
$ws.Range("E25").Value = "  -2.31%  "

Set-TextPrice "D26" "11.82"
$ws.Range("E26").Value = "  -3.66%  "

$ws.Range("E27").Value = "  -6.40%  "

$ws.Range("E28").Value = "  +0.06%  "

Set-TextPrice "D29" "9.98"
$ws.Range("E29").Value = "  -3.02%  "

$ws.Range("E30").Value = "  -2.19%  "

Set-TextPrice "D31" "29.68"
$ws.Range("E31").Value = "  -2.49%  "

$ws.Range("E32").Value = "  -4.45%  "

$ws.Range("E33").Value = "  -3.64%  "

$ws.Range("E34").Value = "  -3.02%  "

$ws.Range("E35").Value = "  -0.07%  "

Set-TextPrice "D36" "3.719.68"
$ws.Range("E36").Value = "  -1.26%  "

$ws.Range("E37").Value = "  -3.20%  "

$ws.Range("E38").Value = "  -8.06%  "

$ws.Range("E39").Value = "  -2.31%  "

$ws.Range("E40").Value = "  -2.45%  "

$ws.Range("E41").Value = "  -3.10%  "

Set-TextPrice "D42" "1.00"
$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("E43").Value = "  +0.01%  "

Set-TextPrice "D44" "43.45"
$ws.Range("E44").Value = "  -1.35%  "

Set-TextPrice "D45" "0.298"
$ws.Range("E45").Value = "  -4.57%  "

Set-TextPrice "D46" "46.83"
$ws.Range("E46").Value = "  +0.82%  "

Set-TextPrice "D47" "8.34"
$ws.Range("E47").Value = "  -3.80%  "

Set-TextPrice "D48" "147.64"
$ws.Range("E48").Value = "  +1.28%  "

$ws.Range("E49").Value = "  -8.74%  "

Set-TextPrice "D50" "389.10"
$ws.Range("E50").Value = "  -3.85%  "

Set-TextPrice "D51" "2.749.44"
$ws.Range("E51").Value = "  +1.74%  "
